# S06/G01: Implement AUTO and MANUAL execution modes per strategy
# Fill in deviations (F), remarks (H), pending work (I) and flip status (G)
# from "pending" to "implemented" for rows 43, 44 and 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# --- Row 43: S06_G01_TB001 ---
$ws.Range("F43").Value = "Strategy execution_mode (AUTO/MANUAL) was already present in the model; we ensured it is consistently exposed through Pydantic schemas and strategy APIs, with MANUAL as the default."
$ws.Range("G43").Value = "implemented"
$ws.Range("H43").Value = "Strategies now carry an execution_mode flag that downstream routing logic and the frontend Settings UI can rely on."
$ws.Range("I43").Value = "Later risk-engine work may add additional per-strategy flags that interact with execution_mode (e.g., risk profiles for AUTO)."

# --- Row 44: S06_G01_TB002 ---
$ws.Range("F44").Value = "TradingView webhook routing now inspects the strategy execution_mode: AUTO strategies (when enabled) create AUTO orders and trigger immediate execution via the existing /api/orders/{id}/execute logic; MANUAL or missing strategies still create WAITING MANUAL queue orders."
$ws.Range("G44").Value = "implemented"
$ws.Range("H44").Value = "AUTO strategies bypass the manual queue and send orders directly to Zerodha using the same execution path as the manual Execute endpoint, including AMO fallback logic."
$ws.Range("I44").Value = "Integrate risk checks into both AUTO and MANUAL routing paths in S06/G02 before broker execution."

# --- Row 45: S06_G01_TF003 ---
$ws.Range("F45").Value = "Settings page now shows an editable Mode column with a MANUAL/AUTO select per strategy, wired to PUT /api/strategies/{id}."
$ws.Range("G45").Value = "implemented"
$ws.Range("H45").Value = "User can toggle strategies between MANUAL and AUTO modes directly from the Settings UI; changes are persisted to the backend."
$ws.Range("I45").Value = "Consider adding per-strategy hints in the UI (e.g., badges or warnings) when AUTO is enabled but broker or risk settings are not fully configured."

$wb.Save()
